# Update "paises" (COVID-19 countries) sheet with a refreshed data pull:
#  - bump the "last updated" timestamp in A1
#  - update case/recovered/death counters for the countries whose ranking
#    shifted between pulls (rows keep their position, but the country that
#    now occupies that rank - and its stats - changes)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 18 de Mayo de 2020 a las 20:05'

$ws.Cells.Item(4, 2).Value = 1535123
$ws.Cells.Item(4, 3).Value = 7459
$ws.Cells.Item(4, 4).Value = 351343
$ws.Cells.Item(4, 5).Value = 1092474

$ws.Cells.Item(6, 2).Value = 278188
$ws.Cells.Item(6, 3).Value = 469
$ws.Cells.Item(6, 5).Value = 54534
$ws.Cells.Item(6, 7).Value = 59
$ws.Cells.Item(6, 8).Value = 27709

$ws.Cells.Item(11, 2).Value = 177213
$ws.Cells.Item(11, 3).Value = 562
$ws.Cells.Item(11, 5).Value = 14517
$ws.Cells.Item(11, 7).Value = 47
$ws.Cells.Item(11, 8).Value = 8096

$ws.Cells.Item(17, 2).Value = 78017
$ws.Cells.Item(17, 3).Value = 1015
$ws.Cells.Item(17, 4).Value = 39127
$ws.Cells.Item(17, 5).Value = 33051
$ws.Cells.Item(17, 7).Value = 57
$ws.Cells.Item(17, 8).Value = 5839

$ws.Cells.Item(32, 1).Value = 'Emiratos Arabes Unidos'
$ws.Cells.Item(32, 2).Value = 24190
$ws.Cells.Item(32, 3).Value = 832
$ws.Cells.Item(32, 4).Value = 9577
$ws.Cells.Item(32, 5).Value = 14389
$ws.Cells.Item(32, 7).Value = 4
$ws.Cells.Item(32, 8).Value = 224

$ws.Cells.Item(33, 1).Value = 'Banglades'
$ws.Cells.Item(33, 2).Value = 23870
$ws.Cells.Item(33, 3).Value = 1602
$ws.Cells.Item(33, 4).Value = 4585
$ws.Cells.Item(33, 5).Value = 18936
$ws.Cells.Item(33, 7).Value = 21
$ws.Cells.Item(33, 8).Value = 349

$ws.Cells.Item(38, 2).Value = 16643
$ws.Cells.Item(38, 3).Value = 26
$ws.Cells.Item(38, 4).Value = 13253
$ws.Cells.Item(38, 5).Value = 3114
$ws.Cells.Item(38, 7).Value = 4
$ws.Cells.Item(38, 8).Value = 276

$ws.Cells.Item(91, 1).Value = 'Consejo Danes para los Refugiados'
$ws.Cells.Item(91, 2).Value = 1538
$ws.Cells.Item(91, 3).Value = 83
$ws.Cells.Item(91, 4).Value = 272
$ws.Cells.Item(91, 5).Value = 1205
$ws.Cells.Item(91, 7).Value = 0
$ws.Cells.Item(91, 8).Value = 61

$ws.Cells.Item(92, 1).Value = 'Republica de Yibuti'
$ws.Cells.Item(92, 2).Value = 1518
$ws.Cells.Item(92, 3).Value = 117
$ws.Cells.Item(92, 4).Value = 1018
$ws.Cells.Item(92, 5).Value = 493
$ws.Cells.Item(92, 7).Value = 3
$ws.Cells.Item(92, 8).Value = 7

$ws.Cells.Item(93, 1).Value = 'Nueva Zelanda'
$ws.Cells.Item(93, 2).Value = 1499
$ws.Cells.Item(93, 3).Value = 0
$ws.Cells.Item(93, 4).Value = 1433
$ws.Cells.Item(93, 5).Value = 45
$ws.Cells.Item(93, 8).Value = 21

$ws.Cells.Item(94, 1).Value = 'Eslovaquia'
$ws.Cells.Item(94, 2).Value = 1495
$ws.Cells.Item(94, 3).Value = 1
$ws.Cells.Item(94, 4).Value = 1185
$ws.Cells.Item(94, 5).Value = 282
$ws.Cells.Item(94, 8).Value = 28

$ws.Cells.Item(95, 1).Value = 'Eslovenia'
$ws.Cells.Item(95, 2).Value = 1466
$ws.Cells.Item(95, 4).Value = 274
$ws.Cells.Item(95, 5).Value = 1088
$ws.Cells.Item(95, 8).Value = 104

$ws.Cells.Item(96, 2).Value = 1455
$ws.Cells.Item(96, 3).Value = 34
$ws.Cells.Item(96, 4).Value = 163
$ws.Cells.Item(96, 5).Value = 1235
$ws.Cells.Item(96, 7).Value = 1
$ws.Cells.Item(96, 8).Value = 57

$ws.Cells.Item(98, 2).Value = 1370
$ws.Cells.Item(98, 3).Value = 28
$ws.Cells.Item(98, 5).Value = 725

$ws.Cells.Item(127, 1).Value = 'Sierra Leona'
$ws.Cells.Item(127, 2).Value = 519
$ws.Cells.Item(127, 3).Value = 14
$ws.Cells.Item(127, 4).Value = 148
$ws.Cells.Item(127, 5).Value = 338
$ws.Cells.Item(127, 7).Value = 1
$ws.Cells.Item(127, 8).Value = 33

$ws.Cells.Item(128, 1).Value = 'Tanzania'
$ws.Cells.Item(128, 2).Value = 509
$ws.Cells.Item(128, 4).Value = 183
$ws.Cells.Item(128, 5).Value = 305
$ws.Cells.Item(128, 8).Value = 21

$ws.Cells.Item(179, 1).Value = 'Zimbabue'
$ws.Cells.Item(179, 2).Value = 46
$ws.Cells.Item(179, 3).Value = 2
$ws.Cells.Item(179, 4).Value = 18
$ws.Cells.Item(179, 5).Value = 24
$ws.Cells.Item(179, 8).Value = 4

$ws.Cells.Item(180, 1).Value = 'Macao'
$ws.Cells.Item(180, 2).Value = 45
$ws.Cells.Item(180, 4).Value = 44
$ws.Cells.Item(180, 5).Value = 1
$ws.Cells.Item(180, 8).Value = 0

$ws.Cells.Item(195, 1).Value = 'Belice'
$ws.Cells.Item(195, 4).Value = 16
$ws.Cells.Item(195, 8).Value = 2

$ws.Cells.Item(196, 1).Value = 'Nueva Caledonia'
$ws.Cells.Item(196, 4).Value = 18
$ws.Cells.Item(196, 8).Value = 0

$ws.Cells.Item(209, 1).Value = 'Groenlandia'

$ws.Cells.Item(210, 1).Value = 'Seychelles'

$ws.Cells.Item(214, 1).Value = 'San Bartolome'

$ws.Cells.Item(215, 1).Value = 'Bonaire, San Eustaquio y Saba'

$ws.Cells.Item(216, 1).Value = 'Sahara Occidental'
